# CALIFORNIA_2015.xlsx cleanup edit
# 1) Rename header columns to snake_case machine-readable names
# 2) Title-case the Spanish grammatical connector words ("de", "del", "el",
#    "la", "los", "las", "y") inside the municipality / state name strings
# 3) Nudge a handful of already-computed percentage values by one ULP
# 4) Drop the trailing footnote rows and shrink the used range back down
#    to the real data (A1:D2363)

$ws = $excel.ActiveWorkbook.ActiveSheet

# --- 1) Header row -----------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2) Title-case connector words in every data row (2 .. 2363) ------
$connectors = @("de", "del", "el", "la", "los", "las", "y")

for ($r = 2; $r -le 2363; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null -and $v -is [string] -and $v.Contains(" ")) {
            $words = $v.Split(" ")
            $changed = $false
            $out = @()
            foreach ($w in $words) {
                if ($connectors -contains $w) {
                    $out += ($w.Substring(0, 1).ToUpper() + $w.Substring(1))
                    $changed = $true
                } else {
                    $out += $w
                }
            }
            if ($changed) {
                $cell.Value = [string]::Join(" ", $out)
            }
        }
    }
}

# --- 3) One-ULP nudges on a few already-rounded percentage cells ------
# (plain decimal literals -- this engine's parser chokes on "E-05" notation)
$ws.Range("D16").Value = 0.00009369411077399604
$ws.Range("D43").Value = 0.00009369411077399604
$ws.Range("D142").Value = 0.009840637340409996
$ws.Range("D188").Value = 0.00009369411077399604
$ws.Range("D299").Value = 0.00009369411077399604
$ws.Range("D333").Value = 0.0009176511437570788
$ws.Range("D422").Value = 0.0009176511437570788
$ws.Range("D600").Value = 0.00009369411077399604
$ws.Range("D616").Value = 0.000939696816880372
$ws.Range("D741").Value = 0.000939696816880372
$ws.Range("D915").Value = 0.0009148954346166672
$ws.Range("D1009").Value = 0.0009259182711783136
$ws.Range("D1108").Value = 0.00009369411077399604
$ws.Range("D1228").Value = 0.00009369411077399604
$ws.Range("D1375").Value = 0.00009369411077399604
$ws.Range("D1514").Value = 0.09288944370499584
$ws.Range("D1636").Value = 0.00009369411077399604
$ws.Range("D1951").Value = 0.00009369411077399604
$ws.Range("D1962").Value = 0.0009644981991440768
$ws.Range("D1979").Value = 0.00009369411077399604
$ws.Range("D2029").Value = 0.00009369411077399604
$ws.Range("D2103").Value = 0.00009369411077399604
$ws.Range("D2169").Value = 0.00009369411077399604
$ws.Range("D2210").Value = 0.00009369411077399604
$ws.Range("D2341").Value = 0.00009369411077399604

# --- 4) Drop the footnote rows (2364 .. 2369) and shrink the dimension -
$ws.Range("A2364:A2369").EntireRow.Delete()
